$d = $word.ActiveDocument
$d.Content.Find.Execute("Square One Standard Employment Offer Letter", $true, $false, $false, $false, $false, $true, 1, $false, "Square One Standard Offer Letter", 2)
